$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Completed"
$ws.Range("E7").Value = "Completed"

$ws.Range("E7").Select()
